$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers I0 and IF in columns I and J, row 1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style (bold, border, centered) from an existing header cell (H1) to I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill column I (I0) with 1 and column J (IF) with the same values as column H, for rows 2-16
for ($r = 2; $r -le 16; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value()
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}
